$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.153.69"
$ws.Range("E2").Value = "  +7.89%  "

$ws.Range("D3").Value = "1.584.67"
$ws.Range("E3").Value = "  +8.05%  "

$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9916"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "297.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3599"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3311"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.761"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.464"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9914"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "

$ws.Range("D17").Value = "1.585.10"
$ws.Range("E17").Value = "  +8.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001054"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06563"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.867"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.01%  "

$ws.Range("D24").Value = "22.149.35"
$ws.Range("E24").Value = "  +7.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.367"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.481"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.76%  "

$ws.Range("D29").Value = "1.755.27"
$ws.Range("E29").Value = "  +7.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.795"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +17.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08068"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.620"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.046"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.228"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02163"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5714"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.750"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.916"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06702"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.62%  "

# Row 39/40: FraxShare and Hedera swap ranking positions with updated values
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.344"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.11%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05948"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.28%  "
